$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59 (pushes old rows 59-99 down to 60-100).
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new weekly record.
$ws.Cells.Item(59, 1).Value = 9
$ws.Cells.Item(59, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(59, 3).Value = "Metropolitana"
$ws.Cells.Item(59, 4).Value = 44658
$ws.Cells.Item(59, 5).Value = 13
$ws.Cells.Item(59, 6).Value = "Fruta"
$ws.Cells.Item(59, 7).Value = 100101
$ws.Cells.Item(59, 8).Value = "Berries"
$ws.Cells.Item(59, 9).Value = 100101004
$ws.Cells.Item(59, 10).Value = "Frambuesa"
$ws.Cells.Item(59, 11).Value = "Sin especificar"
$ws.Cells.Item(59, 12).Value = "Primera"
$ws.Cells.Item(59, 13).Value = 480
$ws.Cells.Item(59, 14).Value = 8000
$ws.Cells.Item(59, 15).Value = 8000
$ws.Cells.Item(59, 16).Value = 8000
$ws.Cells.Item(59, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(59, 18).Value = "Provincia de Linares"
$ws.Cells.Item(59, 19).Value = 4000
$ws.Cells.Item(59, 20).Value = 2
